$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Remove all existing hyperlinks up-front. The row-insert/delete operations
#    below do not reliably shift hyperlink ranges, so we rebuild the full
#    hyperlink collection from scratch at the very end once every row is in
#    its final position.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 1. Insert three new rows (18:20) for M3 hardware (screws/standoffs/nuts).
#    This pushes the old rows 18-42 down to 21-45.
# ---------------------------------------------------------------------------
$ws.Rows("18:20").Insert()

# Row 18: M3 Screws
$ws.Range("A18").Value = 'M3 Screws'
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 'StarTech'
$ws.Range("D18").Value = 'SCREWM3'
$ws.Range("E18").Value = 'Amazon'
$ws.Range("F18").Value = 'B00032Q1J4'
$ws.Range("G18").Value = 0.1034
$ws.Range("H18").Formula = '=B18*G18'
$ws.Range("I18").Value = 'PC Mounting Computer Screws M3x1/4" Long Standoff, 50 Pack'
$ws.Range("J18").Value = 'https://www.amazon.com/dp/B00032Q1J4/_encoding=UTF8?coliid=I3N2TUG8GXM37M&colid=1JXZ7HW0RYXH1&psc=1'

# Row 19: M3 Standoffs (note: vendor part # is a genuine number, kept in Text format)
$ws.Range("A19").Value = 'M3 Standoffs'
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = 'uxcell'
$ws.Range("D19").Value = 'A12092200UX0087'
$ws.Range("E19").Value = 'Amazon'
$ws.Range("F19").Value = 31161816
$ws.Range("F19").NumberFormat = "@"
$ws.Range("G19").Value = 0.1394
$ws.Range("H19").Formula = '=B19*G19'
$ws.Range("I19").Value = '50 Pcs M3 Male x M3 Female Hex PCB Standoffs 12mm Length'
$ws.Range("J19").Value = 'https://www.amazon.com/uxcell%C2%AE-Female-Standoffs-Spacers-Length/dp/B00AH8DEVW/ref=sr_1_11?ie=UTF8&qid=1516654012&sr=8-11&keywords=m3+standoff'

# Row 20: M3 Nuts
$ws.Range("A20").Value = 'M3 Nuts'
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 'Amico'
$ws.Range("D20").Value = ' A13092700UX0955'
$ws.Range("E20").Value = 'Amazon'
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = '700724341402'
$ws.Range("G20").Value = 0.0461
$ws.Range("H20").Formula = '=B20*G20'
$ws.Range("I20").Value = '100Pcs M3 3mm Female Thread Hex Metal Nut Fastener'
$ws.Range("J20").Value = 'https://www.amazon.com/100Pcs-Female-Thread-Fastener-Silver/dp/B00GYS1SXU/ref=sr_1_3?rps=1&ie=UTF8&qid=1516654201&sr=8-3&keywords=m3+nut&refinements=p_85%3A2470955011'

# ---------------------------------------------------------------------------
# 2. Row 16 (Ball Bearings): vendor's part # changes from the manufacturer
#    part # text to a distinct barcode-style value, kept as text.
# ---------------------------------------------------------------------------
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = '692128467461'

# ---------------------------------------------------------------------------
# 3. "Not Included" section (now renamed "Other Costs") at the bottom of the
#    sheet. After the insert above, it currently reads (rows 39-45):
#      39 Not Included / 40 Wiring / 41 Heat Shrink / 42 Solder / 43 PCBs /
#      44 PLA / 45 PCB Standoffs
#    Target layout:
#      39 Other Costs / 40 Wiring / 41 Heat Shrink /
#      42 PCB Fabrication (full line item) / 43 Soldering/Assembly (full line item) /
#      44 PLA / (45,46 blank) / 47 Alternate PCB Assembly Quote
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = 'Other Costs'

# Row 42: replace "Solder" label-only row with a full PCB Fabrication line item
$ws.Range("A42").Value = 'PCB Fabrication'
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 'JLCPCB'
$ws.Range("D42").Value = 'N/A'
$ws.Range("E42").Value = 'N/A'
$ws.Range("F42").Value = 'N/A'
$ws.Range("G42").Value = 0.1909
$ws.Range("H42").Formula = '=B42*G42'
$ws.Range("I42").Value = '*Becomes $0.244/PCB for lead-free RoHS compliance'
$ws.Range("J42").Value = 'https://jlcpcb.com/quote'

# Row 43: replace "PCBs" label-only row with a full Soldering/Assembly line item
$ws.Range("A43").Value = 'Soldering/Assembly'
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 'Screaming Circuits'
$ws.Range("D43").Value = 'N/A'
$ws.Range("E43").Value = 'N/A'
$ws.Range("F43").Value = 'N/A'
$ws.Range("G43").Value = 8.62
$ws.Range("H43").Formula = '=B43*G43'
$ws.Range("I43").Value = '*Just one company quote w/ 15-day lead time'
$ws.Range("J43").Value = 'https://www.screamingcircuits.com/quote'

# Row 45 (old "PCB Standoffs" label, now obsolete since standoffs are a real
# line item in row 19) is removed entirely.
$ws.Rows("45:45").Delete()

# Row 47: new "Alternate PCB Assembly Quote" hyperlink entry (rows 45-46 stay empty)
$ws.Range("A47").Value = 'Alternate PCB Assembly Quote'
$ws.Range("A47").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 4. Rebuild every hyperlink so the collection matches the final row layout.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("J2"), 'https://hobbyking.com/en_us/zippy-compact-850mah-2s-25c-lipo-pack.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J4"), 'https://www.adafruit.com/product/2010') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J5"), 'https://www.amazon.com/Skateboard-Bearings-Double-Shielded-Silver/dp/B002BBGTK6') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J6"), 'https://hobbyking.com/en_us/kingduino-mpu6050-3-axis-gyroscope-3-axis-accelerometer-1.html?countrycode=US&gclid=CjwKCAjw4KvPBRBeEiwAIqCB-bg3Aa4u5bovf1R3RhCnNQGufnsJUiIBR2quZiXUyyfkSDcYw4vmbBoCVYsQAvD_BwE&gclsrc=aw.ds') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J7"), 'https://hobbyking.com/en_us/gemfan-propeller-5x4-black-cw-ccw-2pcs-1.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J8"), 'https://hobbyking.com/en_us/dys-1306-3100kv-bx-series-set-of-two-cw-ccw-motors.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J3"), 'https://hobbyking.com/en_us/turnigy-multistar-10a-v2-esc-with-blheli-and-2a-lbec-2-3s-v.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J17"), 'https://www.amazon.com/LEDMO-Power-Supply-Transformers-Adapter/dp/B01461MOGQ/ref=sr_1_1?ie=UTF8&qid=1516375270&sr=8-1&keywords=B01461MOGQ') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J21"), 'https://hobbyking.com/en_us/brushless-motor-d1306-4000kv-ccw.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J22"), 'https://hobbyking.com/en_us/blheli-s-10a.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J24"), 'https://www.arrow.com/en/products/attiny85-20su/microchip-technology') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J25"), 'https://www.arrow.com/en/products/l78l05abutr/stmicroelectronics') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J28"), 'https://www.arrow.com/en/products/dmn2041l-7/diodes-incorporated') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J29"), 'https://www.arrow.com/en/products/sv03a103aea01r00/murata-manufacturing') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J26"), 'https://www.arrow.com/en/products/c0805c104k3rac7210/kemet-corporation') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J27"), 'https://www.arrow.com/en/products/cl21b334kafnnne/samsung-electro-mechanics') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J23"), 'https://www.arrow.com/en/products/dmn2041l-7/diodes-incorporated') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J30"), 'https://www.arrow.com/en/products/pj-038ah/cui-inc') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J31"), 'https://www.arrow.com/en/products/srb22a2dbbnn/zf-electronics') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J16"), 'https://www.amazon.com/8x22x7-Shielded-Greased-Miniature-Bearings/dp/B00NX3F6F0/ref=pd_sim_328_6?_encoding=UTF8&pd_rd_i=B00NX3F6F0&pd_rd_r=VW2RRR6CAME445GWPFCT&pd_rd_w=EOs5p&pd_rd_wg=tQ5od&psc=1&refRID=VW2RRR6CAME445GWPFCT') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J33"), 'https://www.arrow.com/en/products/1375819-1/te-connectivity') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J32"), 'https://www.arrow.com/en/products/640455-2/te-connectivity') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J34"), 'https://www.arrow.com/en/products/640455-4/te-connectivity') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J35"), 'https://www.arrow.com/en/products/640455-6/te-connectivity') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J36"), 'https://www.arrow.com/en/products/1375820-2/te-connectivity') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J37"), 'https://www.arrow.com/en/products/1375820-6/te-connectivity') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J42"), 'https://jlcpcb.com/quote') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J43"), 'https://www.screamingcircuits.com/quote') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A47"), 'https://www.macrofab.com/') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J18"), 'https://www.amazon.com/dp/B00032Q1J4/_encoding=UTF8?coliid=I3N2TUG8GXM37M&colid=1JXZ7HW0RYXH1&psc=1') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J19"), 'https://www.amazon.com/uxcell%C2%AE-Female-Standoffs-Spacers-Length/dp/B00AH8DEVW/ref=sr_1_11?ie=UTF8&qid=1516654012&sr=8-11&keywords=m3+standoff') | Out-Null
$ws.Hyperlinks.Add($ws.Range("J20"), 'https://www.amazon.com/100Pcs-Female-Thread-Fastener-Silver/dp/B00GYS1SXU/ref=sr_1_3?rps=1&ie=UTF8&qid=1516654201&sr=8-3&keywords=m3+nut&refinements=p_85%3A2470955011') | Out-Null

# ---------------------------------------------------------------------------
# 5. Restore view state as closely as possible.
# ---------------------------------------------------------------------------
$ws.Range("B44").Select()
